$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp string in A1 (22:04 -> 22:34)
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 22:34"

# 2. Swap the country labels for rows 47/48 (Egipto <-> Sudafrica)
#    Row 47 becomes "Sudafrica", row 48 becomes "Egipto"
$ws.Range("A47").Value = "Sudafrica"
$ws.Range("A48").Value = "Egipto"

# 3. Swap the country labels for rows 150/151 (Togo <-> Gibraltar)
#    Row 150 becomes "Gibraltar", row 151 becomes "Togo"
$ws.Range("A150").Value = "Gibraltar"
$ws.Range("A151").Value = "Togo"

# 4. Update the numeric data values that changed between the two data pulls

# Row 9 (Francia)
$ws.Range("B9").Value = 176079
$ws.Range("C9").Value = 1288
$ws.Range("E9").Value = 94067

# Row 17 (Peru)
$ws.Range("D17").Value = 19012
$ws.Range("E17").Value = 37887
$ws.Range("F17").Value = 730

# Row 47 (now Sudafrica)
$ws.Range("B47").Value = 8895
$ws.Range("C47").Value = 663
$ws.Range("D47").Value = 3153
$ws.Range("E47").Value = 5564
$ws.Range("F47").Value = 77
$ws.Range("G47").Value = 17
$ws.Range("H47").Value = 178

# Row 48 (now Egipto)
$ws.Range("B48").Value = 8476
$ws.Range("C48").Value = 495
$ws.Range("D48").Value = 1945
$ws.Range("E48").Value = 6028
$ws.Range("F48").Value = 41
$ws.Range("G48").Value = 21
$ws.Range("H48").Value = 503

# Row 78 (Guinea)
$ws.Range("B78").Value = 2009
$ws.Range("C78").Value = 82
$ws.Range("D78").Value = 663
$ws.Range("E78").Value = 1335

# Row 91 (Republica de Yibuti)
$ws.Range("B91").Value = 1135
$ws.Range("C91").Value = 2
$ws.Range("D91").Value = 824
$ws.Range("E91").Value = 308

# Row 150 (now Gibraltar)
$ws.Range("B150").Value = 146
$ws.Range("C150").Value = 2
$ws.Range("D150").Value = 142
$ws.Range("E150").Value = 4
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 0

# Row 151 (now Togo)
$ws.Range("B151").Value = 145
$ws.Range("C151").Value = 10
$ws.Range("D151").Value = 85
$ws.Range("E151").Value = 50
$ws.Range("G151").Value = 1
$ws.Range("H151").Value = 10
